$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from G1 (bold, bordered, centered) onto the new H1 header cell
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Fill the new "Save" column values for rows 2-5
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 1
$ws.Range("H5").Value = 0
